$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Add available meter types example - set new meter type values for rows 6-10
$ws.Range("D6").Value = "ЭЛ"
$ws.Range("D7").Value = "ЭЛ"
$ws.Range("D8").Value = "ТЕПЛО"
$ws.Range("D9").Value = "ТЕПЛО"
$ws.Range("D10").Value = "ГАЗ"

# Update the active selection as in the authored workbook
$ws.Range("D16").Select()
